$wb = $excel.ActiveWorkbook

# Sheet 1: FEINmismatch
$ws1 = $wb.Worksheets.Item("FEINmismatch")
$ws1.Range("B2").Value = "Thu Jan 25 17:10:25 EST 2024"
$ws1.Range("B3").Value = "Thu Jan 25 17:10:40 EST 2024"
$ws1.Range("B4").Value = "Thu Jan 25 17:10:51 EST 2024"
$ws1.Range("B5").Value = "Thu Jan 25 17:11:04 EST 2024"
$ws1.Range("B6").Value = "Thu Jan 25 17:11:16 EST 2024"
$ws1.Range("B7").Value = "Thu Jan 25 17:11:28 EST 2024"
$ws1.Range("B8").Value = "Thu Jan 25 17:11:40 EST 2024"
$ws1.Range("B9").Value = "Thu Jan 25 17:11:51 EST 2024"
$ws1.Range("B10").Value = "Thu Jan 25 17:12:02 EST 2024"
$ws1.Range("B13").Value = "Thu Jan 25 17:12:14 EST 2024"
$ws1.Range("B14").Value = "Thu Jan 25 17:12:25 EST 2024"
$ws1.Range("B15").Value = "Thu Jan 25 17:12:36 EST 2024"
$ws1.Range("B16").Value = "Thu Jan 25 17:12:48 EST 2024"
$ws1.Range("B17").Value = "Thu Jan 25 17:12:59 EST 2024"
$ws1.Range("B18").Value = "Thu Jan 25 17:13:11 EST 2024"
$ws1.Range("B19").Value = "Thu Jan 25 17:13:27 EST 2024"
$ws1.Range("B20").Value = "Thu Jan 25 17:13:38 EST 2024"
$ws1.Range("B21").Value = "Thu Jan 25 17:13:49 EST 2024"
$ws1.Range("B22").Value = "Thu Jan 25 17:14:01 EST 2024"
$ws1.Range("B23").Value = "Thu Jan 25 17:14:12 EST 2024"
$ws1.Range("B24").Value = "Thu Jan 25 17:14:23 EST 2024"
$ws1.Range("B25").Value = "Thu Jan 25 17:14:34 EST 2024"
$ws1.Range("B26").Value = "Thu Jan 25 17:14:45 EST 2024"
$ws1.Range("B27").Value = "Thu Jan 25 17:14:56 EST 2024"
$ws1.Range("B28").Value = "Thu Jan 25 17:15:08 EST 2024"
$ws1.Range("B29").Value = "Thu Jan 25 17:15:19 EST 2024"
$ws1.Range("B30").Value = "Thu Jan 25 17:15:30 EST 2024"

# Sheet 2: FEINSSNmismatch
$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")
$ws2.Range("B2").Value = "Thu Jan 25 17:15:43 EST 2024"
$ws2.Range("B3").Value = "Thu Jan 25 17:15:54 EST 2024"
$ws2.Range("B4").Value = "Thu Jan 25 17:16:05 EST 2024"
$ws2.Range("B5").Value = "Thu Jan 25 17:16:16 EST 2024"
$ws2.Range("B6").Value = "Thu Jan 25 17:16:27 EST 2024"
$ws2.Range("B7").Value = "Thu Jan 25 17:16:38 EST 2024"
$ws2.Range("B8").Value = "Thu Jan 25 17:16:48 EST 2024"
$ws2.Range("B9").Value = "Thu Jan 25 17:16:59 EST 2024"
$ws2.Range("B14").Value = "Thu Jan 25 17:17:10 EST 2024"
$ws2.Range("B15").Value = "Thu Jan 25 17:17:21 EST 2024"
$ws2.Range("B16").Value = "Thu Jan 25 17:17:32 EST 2024"
$ws2.Range("B17").Value = "Thu Jan 25 17:17:43 EST 2024"
$ws2.Range("B18").Value = "Thu Jan 25 17:17:54 EST 2024"
$ws2.Range("B19").Value = "Thu Jan 25 17:18:05 EST 2024"
